$wb = $excel.ActiveWorkbook

# --- Update the status text from "Ready for handoff" to "In Translation" ---
# This string is shared across the Overview sheet (columns E/F, row 2)
# and the per-locale sheets zh-cn / de-de (column C, row 2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Re-fit the Status columns now that the text is shorter ---
# (the report-generation step auto-sizes these columns to the new content)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5

Write-Output "Updated status text and column widths"
